# Auto-generated Excel COM-interop script
# Applies the numeric corrections described by the commit diff
# to sheets ALC, BSM, CRP, CUL, LTW, WVR.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 239.28572
$ws.Range("I5").Value = 116.666664
$ws.Range("K5").Value = 116.666664
$ws.Range("M5").Value = -1.666663999999997
$ws.Range("H38").Value = 606.5454999999999
$ws.Range("J38").Value = 4500
$ws.Range("L38").Value = 13500
$ws.Range("N38").Value = -14244
$ws.Range("H42").Value = 252.58824
$ws.Range("J42").Value = 305.16666
$ws.Range("L42").Value = 915.4999799999999
$ws.Range("N42").Value = -1375.49998
$ws.Range("H116").Value = 8639.6
$ws.Range("I116").Value = 7733
$ws.Range("K116").Value = 7733
$ws.Range("M116").Value = -4291
$ws.Range("H137").Value = 10423102
$ws.Range("I137").Value = 41668628
$ws.Range("J137").Value = 7926.3887
$ws.Range("K137").Value = 125005884
$ws.Range("L137").Value = 23779.1661
$ws.Range("M137").Value = -125003334
$ws.Range("N137").Value = -28879.1661

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3812.875
$ws.Range("I20").Value = 3707.3333
$ws.Range("K20").Value = 3707.3333
$ws.Range("M20").Value = -3460.3333
$ws.Range("H97").Value = 15999.375
$ws.Range("I97").Value = 10932
$ws.Range("K97").Value = 10932
$ws.Range("M97").Value = -9941
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 55559812
$ws.Range("I31").Value = 125001580
$ws.Range("J31").Value = 6394.2
$ws.Range("K31").Value = 125001580
$ws.Range("L31").Value = 6394.2
$ws.Range("M31").Value = -125001285
$ws.Range("N31").Value = -6984.2
$ws.Range("H34").Value = 55559812
$ws.Range("I34").Value = 125001580
$ws.Range("J34").Value = 6394.2
$ws.Range("K34").Value = 125001580
$ws.Range("L34").Value = 6394.2
$ws.Range("M34").Value = -125001378
$ws.Range("N34").Value = -6798.2
$ws.Range("H54").Value = 35798.6
$ws.Range("J54").Value = 38500
$ws.Range("L54").Value = 38500
$ws.Range("N54").Value = -39816
$ws.Range("H58").Value = 7566.1665
$ws.Range("I58").Value = 5579.4
$ws.Range("J58").Value = 17500
$ws.Range("K58").Value = 5579.4
$ws.Range("L58").Value = 17500
$ws.Range("M58").Value = -5376.4
$ws.Range("N58").Value = -17906
$ws.Range("H134").Value = 8591.639999999999
$ws.Range("I134").Value = 8820.380999999999
$ws.Range("K134").Value = 26461.143
$ws.Range("M134").Value = -23926.143
$ws.Range("H136").Value = 7566.1665
$ws.Range("I136").Value = 5579.4
$ws.Range("J136").Value = 17500
$ws.Range("K136").Value = 16738.2
$ws.Range("L136").Value = 52500
$ws.Range("M136").Value = -14188.2
$ws.Range("N136").Value = -57600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H116").Value = 3750
$ws.Range("I116").Value = 3500
$ws.Range("J116").Value = 4000
$ws.Range("K116").Value = 10500
$ws.Range("L116").Value = 12000
$ws.Range("M116").Value = -7058
$ws.Range("N116").Value = -18884
$ws.Range("H117").Value = 152627.64
$ws.Range("I117").Value = 250
$ws.Range("J117").Value = 167865.4
$ws.Range("K117").Value = 750
$ws.Range("L117").Value = 503596.2
$ws.Range("M117").Value = 2692
$ws.Range("N117").Value = -510480.2
$ws.Range("H121").Value = 23811398
$ws.Range("I121").Value = 376.66666
$ws.Range("K121").Value = 1129.99998
$ws.Range("M121").Value = 180.0000199999999
$ws.Range("H129").Value = 50001350
$ws.Range("I129").Value = 997.1429000000001
$ws.Range("K129").Value = 2991.4287
$ws.Range("M129").Value = 2008.5713
$ws.Range("H131").Value = 15156079
$ws.Range("I131").Value = 33334232
$ws.Range("K131").Value = 100002696
$ws.Range("M131").Value = -99997656

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 5010
$ws.Range("J4").Value = 5010
$ws.Range("L4").Value = 5010
$ws.Range("N4").Value = -5236
$ws.Range("H7").Value = 15679.6
$ws.Range("I7").Value = 11475.125
$ws.Range("K7").Value = 11475.125
$ws.Range("M7").Value = -11363.125
$ws.Range("H22").Value = 3178.7632
$ws.Range("J22").Value = 4390.278
$ws.Range("L22").Value = 4390.278
$ws.Range("N22").Value = -4980.278
$ws.Range("H27").Value = 3178.7632
$ws.Range("J27").Value = 4390.278
$ws.Range("L27").Value = 4390.278
$ws.Range("N27").Value = -4604.278
$ws.Range("H28").Value = 5010
$ws.Range("J28").Value = 5010
$ws.Range("L28").Value = 5010
$ws.Range("N28").Value = -5474
$ws.Range("H30").Value = 7000
$ws.Range("I30").Value = 7000
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 7000
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -6892
$ws.Range("N30").ClearContents()
$ws.Range("H37").Value = 5010
$ws.Range("J37").Value = 5010
$ws.Range("L37").Value = 5010
$ws.Range("N37").Value = -5224
$ws.Range("H46").Value = 4861.8184
$ws.Range("J46").Value = 5429.0713
$ws.Range("L46").Value = 5429.0713
$ws.Range("N46").Value = -5805.0713
$ws.Range("H55").Value = 499.2353
$ws.Range("J55").Value = 547.875
$ws.Range("L55").Value = 547.875
$ws.Range("N55").Value = -893.875
$ws.Range("H64").Value = 386666.66
$ws.Range("J64").Value = 80000
$ws.Range("L64").Value = 80000
$ws.Range("N64").Value = -80450
$ws.Range("H67").Value = 386666.66
$ws.Range("J67").Value = 80000
$ws.Range("L67").Value = 80000
$ws.Range("N67").Value = -81560
$ws.Range("H97").Value = 10332.167
$ws.Range("J97").Value = 10332.167
$ws.Range("L97").Value = 10332.167
$ws.Range("N97").Value = -12314.167
$ws.Range("H103").Value = 20813.428
$ws.Range("J103").Value = 20813.428
$ws.Range("L103").Value = 20813.428
$ws.Range("N103").Value = -23157.428
$ws.Range("H122").Value = 3353.8333
$ws.Range("I122").Value = 3295.0908
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 9885.2724
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -7435.2724
$ws.Range("N122").Value = -16900
$ws.Range("H126").Value = 15679.6
$ws.Range("I126").Value = 11475.125
$ws.Range("K126").Value = 34425.375
$ws.Range("M126").Value = -31955.375
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H130").Value = 47429
$ws.Range("J130").Value = 47429
$ws.Range("L130").Value = 47429
$ws.Range("N130").Value = -57469
$ws.Range("H132").Value = 7623.125
$ws.Range("I132").Value = 6805.385
$ws.Range("J132").Value = 11166.667
$ws.Range("K132").Value = 20416.155
$ws.Range("L132").Value = 33500.001
$ws.Range("M132").Value = -17886.155
$ws.Range("N132").Value = -38560.001
$ws.Range("H133").Value = 49382.625
$ws.Range("J133").Value = 49382.625
$ws.Range("L133").Value = 49382.625
$ws.Range("N133").Value = -54442.625
$ws.Range("H136").Value = 5253.25
$ws.Range("J136").Value = 5558.353
$ws.Range("L136").Value = 16675.059
$ws.Range("N136").Value = -21775.059

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 30405
$ws.Range("I34").Value = 29498
$ws.Range("K34").Value = 29498
$ws.Range("M34").Value = -29295
$ws.Range("H43").Value = 60010
$ws.Range("J43").Value = 85015
$ws.Range("L43").Value = 85015
$ws.Range("N43").Value = -85313
$ws.Range("H52").Value = 9812.75
$ws.Range("I52").Value = 6417.5557
$ws.Range("K52").Value = 6417.5557
$ws.Range("M52").Value = -6191.5557
$ws.Range("H61").Value = 12774.6
$ws.Range("I61").Value = 10706.857
$ws.Range("K61").Value = 10706.857
$ws.Range("M61").Value = -10414.857
$ws.Range("H97").Value = 28916.666
$ws.Range("J97").Value = 28916.666
$ws.Range("L97").Value = 28916.666
$ws.Range("N97").Value = -30898.666
$ws.Range("H113").Value = 522.7646999999999
$ws.Range("I113").Value = 447.18182
$ws.Range("K113").Value = 1341.54546
$ws.Range("M113").Value = 828.45454
